# Edit task "Booking restaurant" (Book / Search form / Logo rows)
#
# Summary of the change (per the commit "Edit task Booking restaurant"):
#   - Row 4  (Book - header row): Build-start date (B4) becomes "bold" styled
#     like the other section headers, and the Build-finish date (C4) is
#     cleared out (task not finished yet).
#   - Row 18 (Book): gets a Build-start date and its Build % raised to 67%.
#   - Row 19 (Search form) and Row 20 (Logo): both get Build-start/finish
#     dates and their Build % set to 100% (done).
#   - The sheet selection/scroll position is moved down to the Book section
#     (cell D18) to reflect where the user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: restyle B4 to the "bold" date format (same look as B10) and
#     clear out C4's finish date ---------------------------------------
$ws.Range("B10").Copy()
$ws.Range("B4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C4").ClearContents()

# --- Row 18 ("Book"): add a start date, bump build % to 0.67 ----------
$ws.Range("B10").Copy()
$ws.Range("B18").PasteSpecial(-4122)  # xlPasteFormats (bold date style)
$ws.Range("B18").Value = 44471
$ws.Range("D18").Value = 0.67

# --- Row 19 ("Search form"): add start/finish dates, build % -> 1 -----
$ws.Range("B5").Copy()
$ws.Range("B19").PasteSpecial(-4122)  # xlPasteFormats (regular date style)
$ws.Range("B19").Value = 44471
$ws.Range("B5").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C19").Value = 44471
$ws.Range("D19").Value = 1

# --- Row 20 ("Logo"): add start/finish dates, build % -> 1 -------------
$ws.Range("B5").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B20").Value = 44471
$ws.Range("B5").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 44471
$ws.Range("D20").Value = 1

$excel.CutCopyMode = $false

# --- Move the on-screen selection / scroll position down to the Book
#     section so the view matches where editing happened -----------------
[void]$ws.Range("D18").Select()
$av = $excel.ActiveWindow
$av.ScrollRow = 16
$av.ScrollColumn = 1
